$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Description - A description of your report (VarBinary(Max)), so you
#    can put it in there if you like" - collapse the spell-check-split runs
#    ("VarBinary" was wrapped in proofErr spellStart/spellEnd markers) back
#    into a single run by replacing the whole sentence with itself.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Description – A description of your report (VarBinary(Max)), so you can put it in there if you like",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Description – A description of your report (VarBinary(Max)), so you can put it in there if you like",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "... margins for the printer and you can change them to fit your
#    printer" -> "... margins for the printer, and you can change them to
#    fit your printer" ; also drop the gramStart/gramEnd proofErr markers
#    that wrapped "printer" and add the comma.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Page Left, Right, Top and Bottom margins. This is the margins for the printer and you can change them to fit your printer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Page Left, Right, Top and Bottom margins. This is the margins for the printer, and you can change them to fit your printer",
    2) | Out-Null

# Re-split "printer," into its own run (matching the target OOXML) by
# toggling Bold on/off over just that word, which forces a run boundary
# without reintroducing any proofErr markers.
$rng = $d.Content
$rng.Find.Execute("printer, and you can change", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$printerStart = $rng.Start
$printerEnd = $printerStart + 8
$printerRun = $d.Range($printerStart, $printerEnd)
$printerRun.Bold = 1
$printerRun.Bold = 0

# ---------------------------------------------------------------------------
# 3. "Fields of Importance" becomes bold (both the run and the paragraph
#    mark pick up bold formatting).
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Fields of Importance*") {
        $p.Range.Bold = 1
        $p.Range.BoldBi = 1
    }
}

# ---------------------------------------------------------------------------
# 4. " To link a data section the linked columns must be available in both
#    canvases. " -> split so "Data Section" is inserted as its own run:
#    " ... in both" + " Data Section" + " canvases. "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " To link a data section the linked columns must be available in both canvases. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " To link a data section the linked columns must be available in both Data Section canvases. ",
    2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute(" To link a data section the linked columns must be available in both Data Section canvases. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sentenceStart = $rng2.Start
$dataSectionStart = $sentenceStart + 68
$dataSectionEnd = $sentenceStart + 81
$dataSectionRun = $d.Range($dataSectionStart, $dataSectionEnd)
$dataSectionRun.Bold = 1
$dataSectionRun.Bold = 0
